$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell as literal text (NumberFormat "@" prevents Excel from
# auto-converting number-like strings), then reset the cell style back to
# "Normal" so no stray style index / quotePrefix is left behind.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "41.741.86"
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -1.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.210.43"
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -1.61%  "
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.38%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "240.76"
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -2.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.623"
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -1.28%  "
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "72.37"
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -4.41%  "
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.20%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.594"
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -4.12%  "
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "41.60"
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -5.64%  "
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0942"
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.73%  "
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "6.92"
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -4.59%  "
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +0.46%  "
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.542.90"
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -1.59%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "14.15"
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -2.88%  "
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.831"
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -3.10%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.217.86"
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -1.96%  "
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "41.650.18"
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -1.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +2.22%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "72.16"
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.13"
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -0.72%  "
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.85"
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +20.38%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "228.49"
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.86%  "
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -9.34%  "
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -0.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.29"
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -1.50%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.26"
$cell.Style = "Normal"

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -2.34%  "
$cell.Style = "Normal"

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.75%  "
$cell.Style = "Normal"

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.56%  "
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "20.33"
$cell.Style = "Normal"

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -1.65%  "
$cell.Style = "Normal"

$cell = $ws.Range("B32")
$cell.NumberFormat = "@"
$cell.Value = "Hedera"
$cell.Style = "Normal"

$cell = $ws.Range("C32")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.0791"
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -4.03%  "
$cell.Style = "Normal"

$cell = $ws.Range("B33")
$cell.NumberFormat = "@"
$cell.Value = "Filecoin"
$cell.Style = "Normal"

$cell = $ws.Range("C33")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.52"
$cell.Style = "Normal"

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +3.54%  "
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "30.06"
$cell.Style = "Normal"

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -2.87%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.124"
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -1.26%  "
$cell.Style = "Normal"

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -11.26%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.20"
$cell.Style = "Normal"

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -7.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0298"
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -5.27%  "
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "13.44"
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -3.61%  "
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.11"
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -3.46%  "
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.60"
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -3.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "63.53"
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -0.73%  "
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -3.19%  "
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "102.51"
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -3.97%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0995"
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -3.43%  "
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -1.62%  "
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -2.49%  "
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -2.82%  "
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -0.97%  "
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.419.33"
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -1.59%  "
$cell.Style = "Normal"

